$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '64.584.19'
$ws.Cells.Item(2, 5).Value = '  -2.52%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.358.60'
$ws.Cells.Item(3, 5).Value = '  -4.53%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '555.29'
$ws.Cells.Item(5, 5).Value = '  -4.77%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '176.22'
$ws.Cells.Item(6, 5).Value = '  -1.88%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.616'
$ws.Cells.Item(7, 5).Value = '  -3.08%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.350.43'
$ws.Cells.Item(8, 5).Value = '  -4.57%  '
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 5).Value = '  -1.94%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.162'
$ws.Cells.Item(11, 5).Value = '  -1.34%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '55.02'
$ws.Cells.Item(12, 5).Value = '  -1.98%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000272'
$ws.Cells.Item(13, 5).Value = '  -3.31%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '9.06'
$ws.Cells.Item(14, 5).Value = '  -2.79%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.911.29'
$ws.Cells.Item(15, 5).Value = '  -4.13%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '18.32'
$ws.Cells.Item(16, 5).Value = '  -0.64%  '
$ws.Cells.Item(17, 5).Value = '  -2.71%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.355.51'
$ws.Cells.Item(18, 5).Value = '  -4.57%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.81'
$ws.Cells.Item(19, 5).Value = '  -2.21%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '64.565.12'
$ws.Cells.Item(20, 5).Value = '  -2.61%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.979'
$ws.Cells.Item(21, 5).Value = '  -3.54%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '436.74'
$ws.Cells.Item(22, 5).Value = '  +4.99%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.96'
$ws.Cells.Item(23, 5).Value = '  +12.49%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '4.08'
$ws.Cells.Item(24, 5).Value = '  -5.08%  '
$ws.Cells.Item(25, 5).Value = '  -0.97%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '13.35'
$ws.Cells.Item(26, 5).Value = '  -1.62%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.79'
$ws.Cells.Item(27, 5).Value = '  -2.95%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.83'
$ws.Cells.Item(28, 5).Value = '  -1.38%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.77'
$ws.Cells.Item(29, 5).Value = '  -4.79%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '29.73'
$ws.Cells.Item(30, 5).Value = '  -2.36%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '6.61'
$ws.Cells.Item(31, 5).Value = '  -0.93%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '11.47'
$ws.Cells.Item(32, 5).Value = '  -2.76%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '579.80'
$ws.Cells.Item(33, 5).Value = '  -4.23%  '
$ws.Cells.Item(34, 5).Value = '  -3.27%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '58.53'
$ws.Cells.Item(35, 5).Value = '  -3.64%  '
$ws.Cells.Item(36, 5).Value = '  +0.06%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.142'
$ws.Cells.Item(37, 5).Value = '  -7.82%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.52'
$ws.Cells.Item(38, 5).Value = '  -4.10%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '35.73'
$ws.Cells.Item(39, 5).Value = '  -3.41%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0₃0755'
$ws.Cells.Item(40, 5).Value = '  -6.11%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.368'
$ws.Cells.Item(41, 5).Value = '  -4.93%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.108.69'
$ws.Cells.Item(42, 5).Value = '  -4.24%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  +0.12%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.82'
$ws.Cells.Item(44, 5).Value = '  -5.46%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.23'
$ws.Cells.Item(45, 5).Value = '  -3.32%  '
$ws.Cells.Item(46, 5).Value = '  -2.87%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.46'
$ws.Cells.Item(47, 5).Value = '  -4.25%  '
$ws.Cells.Item(48, 5).Value = '  -2.71%  '
$ws.Cells.Item(49, 5).Value = '  -3.88%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '8.31'
$ws.Cells.Item(50, 5).Value = '  -4.20%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '135.18'
$ws.Cells.Item(51, 5).Value = '  -2.44%  '
